$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44832
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13500
$ws.Range("P2").Value = 1038

# Row 3
$ws.Range("D3").Value = 45141
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("L3").Value = 17000
$ws.Range("M3").Value = 16550
$ws.Range("P3").Value = 1273

# Row 4
$ws.Range("D4").Value = 44469
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("N4").Value = '$/caja 13 kilos'
$ws.Range("P4").Value = 1038
$ws.Range("Q4").Value = 13

# Row 5
$ws.Range("D5").Value = 44320
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 19000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 19500
$ws.Range("P5").Value = 1500

# Row 7
$ws.Range("D7").Value = 45154
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("P7").Value = 1346

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 23000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 23500
$ws.Range("P8").Value = 1808

# Row 9
$ws.Range("D9").Value = 45049
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("P9").Value = 1038

# Row 10
$ws.Range("D10").Value = 44984
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 16500
$ws.Range("P10").Value = 1269

# Row 11
$ws.Range("D11").Value = 44406
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 1346

# Row 12
$ws.Range("D12").Value = 44910
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("P12").Value = 1115

# Row 13
$ws.Range("D13").Value = 44397
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12750
$ws.Range("P13").Value = 981

# Row 14
$ws.Range("D14").Value = 44379
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 12667
$ws.Range("P14").Value = 974

# Row 15
$ws.Range("D15").Value = 44389
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12500
$ws.Range("P15").Value = 962

# Row 16
$ws.Range("D16").Value = 44890
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("P16").Value = 1115

# Row 17
$ws.Range("D17").Value = 44580
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11500
$ws.Range("P17").Value = 885

# Row 18
$ws.Range("D18").Value = 45100
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 1192

# Row 19
$ws.Range("D19").Value = 45028
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("P19").Value = 1115

# Row 20
$ws.Range("D20").Value = 44893
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13444
$ws.Range("P20").Value = 1034

# Row 21
$ws.Range("D21").Value = 44764
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 13000
$ws.Range("M21").Value = 12500
$ws.Range("P21").Value = 962

# Row 22
$ws.Range("D22").Value = 44972
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17429
$ws.Range("N22").Value = '$/caja 15 kilos'
$ws.Range("P22").Value = 1162
$ws.Range("Q22").Value = 15

# Row 23
$ws.Range("D23").Value = 45092
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 13500
$ws.Range("P23").Value = 1038

# Row 24
$ws.Range("D24").Value = 44943
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 350
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14429
$ws.Range("P24").Value = 1110

# Row 25
$ws.Range("D25").Value = 44592
$ws.Range("J25").Value = 120

# Row 26
$ws.Range("D26").Value = 44914
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14400
$ws.Range("P26").Value = 1108

# Row 27
$ws.Range("D27").Value = 44918
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 12000
$ws.Range("L27").Value = 13000
$ws.Range("M27").Value = 12750
$ws.Range("P27").Value = 981

# Row 28
$ws.Range("D28").Value = 45096
$ws.Range("J28").Value = 750
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14600
$ws.Range("P28").Value = 1123

# Row 29
$ws.Range("D29").Value = 44616
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 19000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 19500
$ws.Range("P29").Value = 1500

# Row 30
$ws.Range("D30").Value = 44988
$ws.Range("J30").Value = 750
$ws.Range("K30").Value = 17000
$ws.Range("L30").Value = 18000
$ws.Range("M30").Value = 17400
$ws.Range("P30").Value = 1338

# Row 31
$ws.Range("D31").Value = 45140
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 16000
$ws.Range("L31").Value = 17000
$ws.Range("M31").Value = 16500
$ws.Range("P31").Value = 1269

# Row 32
$ws.Range("D32").Value = 44855
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("P32").Value = 769
